# This workbook recomputes NATMI ligand-receptor edge-weight statistics for
# the Tgfb1-Tgfbr3 pair using updated per-cluster TPM values. Ligand (G/H/I/J)
# and receptor (M/N/O/P) average/total expression & derived-specificity values
# change per Sending/Target cluster, and the resulting edge weights and their
# derived specificities (Q/R/S/T) are updated for every Sending x Target row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 54.53585066666667
$ws.Range("H2").Value = 163.607552
$ws.Range("I2").Value = 0.3031388658437607
$ws.Range("J2").Value = 0.3031388658437607
$ws.Range("M2").Value = 42.09975866666667
$ws.Range("N2").Value = 126.299276
$ws.Range("O2").Value = 0.3315552933456474
$ws.Range("P2").Value = 0.3315552933456474
$ws.Range("Q2").Value = 2295.946151748039
$ws.Range("R2").Value = 20663.51536573235
$ws.Range("S2").Value = 0.1005072955892949
$ws.Range("T2").Value = 0.1005072955892949

# Row 3
$ws.Range("G3").Value = 54.53585066666667
$ws.Range("H3").Value = 163.607552
$ws.Range("I3").Value = 0.3031388658437607
$ws.Range("J3").Value = 0.3031388658437607
$ws.Range("O3").Value = 0.4502223747274475
$ws.Range("P3").Value = 0.4502223747274475
$ws.Range("Q3").Value = 3117.689113799567
$ws.Range("R3").Value = 28059.2020241961
$ws.Range("S3").Value = 0.136479900052363
$ws.Range("T3").Value = 0.1364799000523631

# Row 4
$ws.Range("G4").Value = 54.53585066666667
$ws.Range("H4").Value = 163.607552
$ws.Range("I4").Value = 0.3031388658437607
$ws.Range("J4").Value = 0.3031388658437607
$ws.Range("M4").Value = 27.596267
$ws.Range("N4").Value = 82.78880100000001
$ws.Range("O4").Value = 0.2173335118824389
$ws.Range("P4").Value = 0.2173335118824389
$ws.Range("Q4").Value = 1504.985896069461
$ws.Range("R4").Value = 13544.87306462515
$ws.Range("S4").Value = 0.06588223430188399
$ws.Range("T4").Value = 0.06588223430188402

# Row 5
$ws.Range("G5").Value = 54.53585066666667
$ws.Range("H5").Value = 163.607552
$ws.Range("I5").Value = 0.3031388658437607
$ws.Range("J5").Value = 0.3031388658437607
$ws.Range("M5").Value = 0.1128593333333333
$ws.Range("N5").Value = 0.338578
$ws.Range("O5").Value = 0.0008888200444663087
$ws.Range("P5").Value = 0.0008888200444663087
$ws.Range("Q5").Value = 6.154879749006222
$ws.Range("R5").Value = 55.393917741056
$ws.Range("S5").Value = 0.0002694359002187177
$ws.Range("T5").Value = 0.0002694359002187178

# Row 6
$ws.Range("I6").Value = 0.1026363515063155
$ws.Range("J6").Value = 0.1026363515063155
$ws.Range("M6").Value = 42.09975866666667
$ws.Range("N6").Value = 126.299276
$ws.Range("O6").Value = 0.3315552933456474
$ws.Range("P6").Value = 0.3315552933456474
$ws.Range("Q6").Value = 777.3583753916862
$ws.Range("R6").Value = 6996.225378525175
$ws.Range("S6").Value = 0.03402962563160341
$ws.Range("T6").Value = 0.03402962563160342

# Row 7
$ws.Range("I7").Value = 0.1026363515063155
$ws.Range("J7").Value = 0.1026363515063155
$ws.Range("O7").Value = 0.4502223747274475
$ws.Range("P7").Value = 0.4502223747274475
$ws.Range("S7").Value = 0.04620918190853439
$ws.Range("T7").Value = 0.0462091819085344

# Row 8
$ws.Range("I8").Value = 0.1026363515063155
$ws.Range("J8").Value = 0.1026363515063155
$ws.Range("M8").Value = 27.596267
$ws.Range("N8").Value = 82.78880100000001
$ws.Range("O8").Value = 0.2173335118824389
$ws.Range("P8").Value = 0.2173335118824389
$ws.Range("Q8").Value = 509.5561105669806
$ws.Range("R8").Value = 4586.004995102826
$ws.Range("S8").Value = 0.02230631871966799
$ws.Range("T8").Value = 0.02230631871966799

# Row 9
$ws.Range("I9").Value = 0.1026363515063155
$ws.Range("J9").Value = 0.1026363515063155
$ws.Range("M9").Value = 0.1128593333333333
$ws.Range("N9").Value = 0.338578
$ws.Range("O9").Value = 0.0008888200444663087
$ws.Range("P9").Value = 0.0008888200444663087
$ws.Range("Q9").Value = 2.083910948336444
$ws.Range("R9").Value = 18.755198535028
$ws.Range("S9").Value = 0.00009122524650970303
$ws.Range("T9").Value = 0.00009122524650970304

# Row 10
$ws.Range("G10").Value = 12.55635966666667
$ws.Range("H10").Value = 37.669079
$ws.Range("I10").Value = 0.06979483370938171
$ws.Range("J10").Value = 0.06979483370938172
$ws.Range("M10").Value = 42.09975866666667
$ws.Range("N10").Value = 126.299276
$ws.Range("O10").Value = 0.3315552933456474
$ws.Range("P10").Value = 0.3315552933456474
$ws.Range("Q10").Value = 528.6197116985337
$ws.Range("R10").Value = 4757.577405286804
$ws.Range("S10").Value = 0.02314084656452473
$ws.Range("T10").Value = 0.02314084656452474

# Row 11
$ws.Range("G11").Value = 12.55635966666667
$ws.Range("H11").Value = 37.669079
$ws.Range("I11").Value = 0.06979483370938171
$ws.Range("J11").Value = 0.06979483370938172
$ws.Range("O11").Value = 0.4502223747274475
$ws.Range("P11").Value = 0.4502223747274475
$ws.Range("Q11").Value = 717.8181941451935
$ws.Range("R11").Value = 6460.363747306742
$ws.Range("S11").Value = 0.03142319577634513
$ws.Range("T11").Value = 0.03142319577634514

# Row 12
$ws.Range("G12").Value = 12.55635966666667
$ws.Range("H12").Value = 37.669079
$ws.Range("I12").Value = 0.06979483370938171
$ws.Range("J12").Value = 0.06979483370938172
$ws.Range("M12").Value = 27.596267
$ws.Range("N12").Value = 82.78880100000001
$ws.Range("O12").Value = 0.2173335118824389
$ws.Range("P12").Value = 0.2173335118824389
$ws.Range("Q12").Value = 346.5086539093643
$ws.Range("R12").Value = 3118.577885184279
$ws.Range("S12").Value = 0.01516875632131075
$ws.Range("T12").Value = 0.01516875632131076

# Row 13
$ws.Range("G13").Value = 12.55635966666667
$ws.Range("H13").Value = 37.669079
$ws.Range("I13").Value = 0.06979483370938171
$ws.Range("J13").Value = 0.06979483370938172
$ws.Range("M13").Value = 0.1128593333333333
$ws.Range("N13").Value = 0.338578
$ws.Range("O13").Value = 0.0008888200444663087
$ws.Range("P13").Value = 0.0008888200444663087
$ws.Range("Q13").Value = 1.417102381073555
$ws.Range("R13").Value = 12.753921429662
$ws.Range("S13").Value = 0.00006203504720109127
$ws.Range("T13").Value = 0.00006203504720109129

# Row 14
$ws.Range("G14").Value = 94.34696966666667
$ws.Range("H14").Value = 283.040909
$ws.Range("I14").Value = 0.524429948940542
$ws.Range("J14").Value = 0.5244299489405421
$ws.Range("M14").Value = 42.09975866666667
$ws.Range("N14").Value = 126.299276
$ws.Range("O14").Value = 0.3315552933456474
$ws.Range("P14").Value = 0.3315552933456474
$ws.Range("Q14").Value = 3971.984653897987
$ws.Range("R14").Value = 35747.86188508188
$ws.Range("S14").Value = 0.1738775255602243
$ws.Range("T14").Value = 0.1738775255602243

# Row 15
$ws.Range("G15").Value = 94.34696966666667
$ws.Range("H15").Value = 283.040909
$ws.Range("I15").Value = 0.524429948940542
$ws.Range("J15").Value = 0.5244299489405421
$ws.Range("O15").Value = 0.4502223747274475
$ws.Range("P15").Value = 0.4502223747274475
$ws.Range("Q15").Value = 5393.59919491512
$ws.Range("R15").Value = 48542.39275423609
$ws.Range("S15").Value = 0.2361100969902048
$ws.Range("T15").Value = 0.2361100969902049

# Row 16
$ws.Range("G16").Value = 94.34696966666667
$ws.Range("H16").Value = 283.040909
$ws.Range("I16").Value = 0.524429948940542
$ws.Range("J16").Value = 0.5244299489405421
$ws.Range("M16").Value = 27.596267
$ws.Range("N16").Value = 82.78880100000001
$ws.Range("O16").Value = 0.2173335118824389
$ws.Range("P16").Value = 0.2173335118824389
$ws.Range("Q16").Value = 2603.624165562234
$ws.Range("R16").Value = 23432.61749006011
$ws.Range("S16").Value = 0.1139762025395761
$ws.Range("T16").Value = 0.1139762025395761

# Row 17
$ws.Range("G17").Value = 94.34696966666667
$ws.Range("H17").Value = 283.040909
$ws.Range("I17").Value = 0.524429948940542
$ws.Range("J17").Value = 0.5244299489405421
$ws.Range("M17").Value = 0.1128593333333333
$ws.Range("N17").Value = 0.338578
$ws.Range("O17").Value = 0.0008888200444663087
$ws.Range("P17").Value = 0.0008888200444663087
$ws.Range("Q17").Value = 10.64793609860022
$ws.Range("R17").Value = 95.83142488740199
$ws.Range("S17").Value = 0.0004661238505367965
$ws.Range("T17").Value = 0.0004661238505367966

